# Applies the cryptos.xlsx data refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = '@'
$r.Value = '71.423.09'
$r.ClearFormats()
$ws.Range('E2').Value = '  -1.43%  '
$r = $ws.Range('D3')
$r.NumberFormat = '@'
$r.Value = '3.978.28'
$r.ClearFormats()
$ws.Range('E3').Value = '  -1.93%  '
$ws.Range('E4').Value = '  -0.10%  '
$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '541.34'
$r.ClearFormats()
$ws.Range('E5').Value = '  +4.27%  '
$r = $ws.Range('D6')
$r.NumberFormat = '@'
$r.Value = '148.94'
$r.ClearFormats()
$ws.Range('E6').Value = '  +1.12%  '
$r = $ws.Range('D7')
$r.NumberFormat = '@'
$r.Value = '3.969.35'
$r.ClearFormats()
$ws.Range('E7').Value = '  -1.93%  '
$ws.Range('E8').Value = '  -6.26%  '
$r = $ws.Range('D10')
$r.NumberFormat = '@'
$r.Value = '0.743'
$r.ClearFormats()
$ws.Range('E10').Value = '  -3.77%  '
$ws.Range('E11').Value = '  -5.07%  '
$r = $ws.Range('D12')
$r.NumberFormat = '@'
$r.Value = '56.29'
$r.ClearFormats()
$ws.Range('E12').Value = '  +17.98%  '
$r = $ws.Range('D13')
$r.NumberFormat = '@'
$r.Value = '0.0000319'
$r.ClearFormats()
$ws.Range('E13').Value = '  -2.78%  '
$r = $ws.Range('D14')
$r.NumberFormat = '@'
$r.Value = '10.75'
$r.ClearFormats()
$ws.Range('E14').Value = '  -3.56%  '
$r = $ws.Range('D15')
$r.NumberFormat = '@'
$r.Value = '4.612.63'
$r.ClearFormats()
$ws.Range('E15').Value = '  -1.95%  '
$r = $ws.Range('D16')
$r.NumberFormat = '@'
$r.Value = '3.986.46'
$r.ClearFormats()
$ws.Range('E16').Value = '  -2.20%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$r = $ws.Range('D17')
$r.NumberFormat = '@'
$r.Value = '20.69'
$r.ClearFormats()
$ws.Range('E17').Value = '  -2.40%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$r = $ws.Range('D18')
$r.NumberFormat = '@'
$r.Value = '13.96'
$r.ClearFormats()
$ws.Range('E18').Value = '  -1.07%  '
$ws.Range('E19').Value = '  -1.14%  '
$ws.Range('E20').Value = '  -2.69%  '
$r = $ws.Range('D21')
$r.NumberFormat = '@'
$r.Value = '71.362.90'
$r.ClearFormats()
$ws.Range('E21').Value = '  -1.50%  '
$r = $ws.Range('D22')
$r.NumberFormat = '@'
$r.Value = '428.87'
$r.ClearFormats()
$ws.Range('E22').Value = '  -3.28%  '
$ws.Range('B23').Value = 'ImmutableX'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$r = $ws.Range('D23')
$r.NumberFormat = '@'
$r.Value = '3.60'
$r.ClearFormats()
$ws.Range('E23').Value = '  +1.03%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$r = $ws.Range('D24')
$r.NumberFormat = '@'
$r.Value = '97.61'
$r.ClearFormats()
$ws.Range('E24').Value = '  -6.82%  '
$ws.Range('E25').Value = '  +5.24%  '
$r = $ws.Range('D26')
$r.NumberFormat = '@'
$r.Value = '14.50'
$r.ClearFormats()
$ws.Range('E26').Value = '  -2.29%  '
$r = $ws.Range('D27')
$r.NumberFormat = '@'
$r.Value = '11.54'
$r.ClearFormats()
$ws.Range('E27').Value = '  +0.67%  '
$r = $ws.Range('D28')
$r.NumberFormat = '@'
$r.Value = '10.76'
$r.ClearFormats()
$ws.Range('E28').Value = '  -2.40%  '
$ws.Range('E29').Value = '  +14.51%  '
$ws.Range('E30').Value = '  +1.83%  '
$r = $ws.Range('D31')
$r.NumberFormat = '@'
$r.Value = '36.68'
$r.ClearFormats()
$ws.Range('E31').Value = '  -3.02%  '
$r = $ws.Range('D32')
$r.NumberFormat = '@'
$r.Value = '7.75'
$r.ClearFormats()
$ws.Range('E32').Value = '  +13.71%  '
$r = $ws.Range('D33')
$r.NumberFormat = '@'
$r.Value = '51.50'
$r.ClearFormats()
$ws.Range('E33').Value = '  +21.38%  '
$r = $ws.Range('D34')
$r.NumberFormat = '@'
$r.Value = '694.83'
$r.ClearFormats()
$ws.Range('E34').Value = '  +1.81%  '
$r = $ws.Range('D35')
$r.NumberFormat = '@'
$r.Value = '13.47'
$r.ClearFormats()
$ws.Range('E35').Value = '  -1.57%  '
$ws.Range('E36').Value = '  +0.79%  '
$r = $ws.Range('D37')
$r.NumberFormat = '@'
$r.Value = '65.61'
$r.ClearFormats()
$ws.Range('E37').Value = '  -2.13%  '
$r = $ws.Range('D38')
$r.NumberFormat = '@'
$r.Value = '0.439'
$r.ClearFormats()
$ws.Range('E38').Value = '  +2.51%  '
$r = $ws.Range('D39')
$r.NumberFormat = '@'
$r.Value = '0.0₃0826'
$r.ClearFormats()
$ws.Range('E39').Value = '  -4.52%  '
$ws.Range('E40').Value = '  +0.60%  '
$r = $ws.Range('D41')
$r.NumberFormat = '@'
$r.Value = '3.42'
$r.ClearFormats()
$ws.Range('E41').Value = '  -2.80%  '
$ws.Range('E42').Value = '  +0.47%  '
$ws.Range('E43').Value = '  +0.11%  '
$r = $ws.Range('D44')
$r.NumberFormat = '@'
$r.Value = '3.29'
$r.ClearFormats()
$ws.Range('E44').Value = '  +0.53%  '
$ws.Range('E45').Value = '  -2.37%  '
$ws.Range('E46').Value = '  -6.18%  '
$r = $ws.Range('D47')
$r.NumberFormat = '@'
$r.Value = '2.72'
$r.ClearFormats()
$ws.Range('E47').Value = '  +1.33%  '
$r = $ws.Range('D48')
$r.NumberFormat = '@'
$r.Value = '9.75'
$r.ClearFormats()
$ws.Range('E48').Value = '  +6.19%  '
$r = $ws.Range('D49')
$r.NumberFormat = '@'
$r.Value = '3.39'
$r.ClearFormats()
$ws.Range('E49').Value = '  -3.15%  '
$r = $ws.Range('D50')
$r.NumberFormat = '@'
$r.Value = '3.00'
$r.ClearFormats()
$ws.Range('E50').Value = '  -2.02%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$r = $ws.Range('D51')
$r.NumberFormat = '@'
$r.Value = '2.758.50'
$r.ClearFormats()
$ws.Range('E51').Value = '  +5.74%  '
